# This script applies the commit "Fruta / hortaliza, semanal": two new
# weekly price records are inserted into the Acelga (Macroferia Regional
# de Talca) dataset. The new records are inserted as new rows 108 and 184
# (pushing all subsequent data rows down), which matches the shift seen
# between the pre- and post-commit data (everything from the former row
# 108 onward moves down by one row after the first insertion, and
# everything from the former row 183 onward moves down by one more row
# after the second insertion).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the two new rows first (this shifts all the existing data rows
# down and automatically keeps their D/J/K/L/M/P (and every other)
# values correct relative to their new row numbers).
$ws.Rows.Item(108).Insert()
$ws.Rows.Item(184).Insert()

# Populate the first new row (new row 108) with its full record.
$ws.Cells.Item(108, 1).Value = 5
$ws.Cells.Item(108, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(108, 3).Value = "Maule"
$ws.Cells.Item(108, 4).Value = 44664
$ws.Cells.Item(108, 5).Value = 7
$ws.Cells.Item(108, 6).Value = 100112009
$ws.Cells.Item(108, 7).Value = "Acelga"
$ws.Cells.Item(108, 8).Value = "Sin especificar"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 400
$ws.Cells.Item(108, 11).Value = 3500
$ws.Cells.Item(108, 12).Value = 3500
$ws.Cells.Item(108, 13).Value = 3500
$ws.Cells.Item(108, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(108, 15).Value = "Región del Maule"
$ws.Cells.Item(108, 16).Value = 875
$ws.Cells.Item(108, 17).Value = 4
$ws.Cells.Item(108, 18).Value = "Hortaliza"
$ws.Cells.Item(108, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the second new row (new row 184) with its full record.
$ws.Cells.Item(184, 1).Value = 5
$ws.Cells.Item(184, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(184, 3).Value = "Maule"
$ws.Cells.Item(184, 4).Value = 44663
$ws.Cells.Item(184, 5).Value = 7
$ws.Cells.Item(184, 6).Value = 100112009
$ws.Cells.Item(184, 7).Value = "Acelga"
$ws.Cells.Item(184, 8).Value = "Sin especificar"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 500
$ws.Cells.Item(184, 11).Value = 3500
$ws.Cells.Item(184, 12).Value = 3500
$ws.Cells.Item(184, 13).Value = 3500
$ws.Cells.Item(184, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(184, 15).Value = "Región del Maule"
$ws.Cells.Item(184, 16).Value = 875
$ws.Cells.Item(184, 17).Value = 4
$ws.Cells.Item(184, 18).Value = "Hortaliza"
$ws.Cells.Item(184, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Host "Inserted two new weekly price rows (108 and 184)."
